# A new price-report row was inserted into the daily Mango price log at
# sheet row 539. Every existing row from 539 down to 651 shifts down by one
# (i.e. to 540-652), and the newly opened row 539 receives its own
# Fecha / Volumen / Precio minimo / Precio maximo / Precio promedio
# ponderado / Origen / Precio $/Kg values, while the "template" columns
# that are constant for this whole report (Mercado ID, Mercado, Region,
# Codreg, Tipo, Producto ID, Producto, Categoria ID, Categoria, Variedad,
# Calidad, Unidad de comercializacion, Kg / unidad) keep the same pattern
# they always have in this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert one blank row above the current row 539; this pushes rows
# 539..651 down to 540..652 and extends the sheet dimension to A1:T652.
$ws.Rows.Item(539).Insert()

# Fill in the freshly inserted row 539 with the new record's data.
$ws.Cells.Item(539, 1).Value = 10
$ws.Cells.Item(539, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(539, 3).Value = "La Araucanía"
$ws.Cells.Item(539, 4).Value = 45209
$ws.Cells.Item(539, 5).Value = 9
$ws.Cells.Item(539, 6).Value = "Fruta"
$ws.Cells.Item(539, 7).Value = 100108
$ws.Cells.Item(539, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(539, 9).Value = 100108002
$ws.Cells.Item(539, 10).Value = "Mango"
$ws.Cells.Item(539, 11).Value = "Sin especificar"
$ws.Cells.Item(539, 12).Value = "Primera"
$ws.Cells.Item(539, 13).Value = 250
$ws.Cells.Item(539, 14).Value = 10000
$ws.Cells.Item(539, 15).Value = 10000
$ws.Cells.Item(539, 16).Value = 10000
$ws.Cells.Item(539, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(539, 18).Value = "Brasil"
$ws.Cells.Item(539, 19).Value = 2500
$ws.Cells.Item(539, 20).Value = 4
